$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts existing rows 2-5 down to 3-6)
$ws.Rows.Item(2).Insert()

# Fill in the new row 2 with the TEST_BARBARA shipment data
$ws.Cells.Item(2, 1).Value = "TEST_BARBARA"
$ws.Cells.Item(2, 2).Value = "dpd-nl"
$ws.Cells.Item(2, 3).Value = "DPD Netherlands"
$ws.Cells.Item(2, 4).Value = "Delivered"
$ws.Cells.Item(2, 5).Value = ""
$ws.Cells.Item(2, 6).Value = ""
$ws.Cells.Item(2, 7).Value = "2026-02-09T00:30:00"
$ws.Cells.Item(2, 8).Value = ""
$ws.Cells.Item(2, 9).Value = "2026-02-09T07:16:12+00:00"
$ws.Cells.Item(2, 10).Value = "TEST_BARBARA"
$ws.Cells.Item(2, 11).Value = "{}"
